# Fruta / hortaliza, semanal
# Inserts this week's two new price records (row 3 and row 4) for
# "Feria Lagunitas de Puerto Montt - Chirimoya", pushing the existing
# historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 3 (rows 3 and 4), shifting all
# the existing data (previously rows 3..16) down to rows 5..18.
$ws.Range("A3:A4").EntireRow.Insert()

# New row 3
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44530
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("Q3").Value = "$/bandeja 8 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 2438
$ws.Range("T3").Value = 8

# New row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44530
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 2000
$ws.Range("T4").Value = 8
